# Updated mapping from SEDOL to ISIN
# - Switch calculation mode back to Automatic (drops calcMode="manual")
# - Append four new lookup rows (65-68) with their SEDOL / bond-name / ISIN values,
#   copying the existing formatting used by the table so new cells match neighbours
# - Move the active selection to C68 to match the refreshed sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Workbook was left in manual calc mode; switch back to automatic.
$excel.Calculation = -4105

# New SEDOL -> Name -> ISIN rows appended to the lookup table.
# Values are written column-by-column for rows 65-67 (matching how the shared
# strings table ended up ordered), then row 68 is filled in row-by-row.
$ws.Range("A65").Value = "BDDWMY1"
$ws.Range("A66").Value = "BDTYZ31"
$ws.Range("A67").Value = "BF07339"

$ws.Range("B65").Value = "SOPOWZ 3 1/2 05/08/27"
$ws.Range("B66").Value = "LOGPH 5 3/4 01/03/22"
$ws.Range("B67").Value = "CHJMAO 4 PERP"

$ws.Range("C65").Value = "USG2120QAC09"
$ws.Range("C66").Value = "XS1541978851"
$ws.Range("C67").Value = "XS1673588452"

$ws.Range("A68").Value = "BYW5T25"
$ws.Range("B68").Value = "KAISAG 8 1/2 06/30/22"
$ws.Range("C68").Value = "XS1627597955"

# Column A carries the existing "code" formatting for all four new rows; only the
# last row (68) had B/C formatted to match as well (mirrors the source edit).
$ws.Range("A65").Style = $ws.Range("A64").Style
$ws.Range("A66").Style = $ws.Range("A64").Style
$ws.Range("A67").Style = $ws.Range("A64").Style
$ws.Range("A68").Style = $ws.Range("A64").Style
$ws.Range("B68").Style = $ws.Range("B64").Style
$ws.Range("C68").Style = $ws.Range("C64").Style

# Match the active-cell selection saved in the sheet view.
$ws.Range("C68").Select() | Out-Null
